$d = $word.ActiveDocument

$d.Content.Find.Execute("79×75=5925", $true, $false, $false, $false, $false, $true, 1, $false, "63×20=1260", 2) | Out-Null
$d.Content.Find.Execute("38×56=2128", $true, $false, $false, $false, $false, $true, 1, $false, "17×77=1309", 2) | Out-Null
$d.Content.Find.Execute("75×25=1875", $true, $false, $false, $false, $false, $true, 1, $false, "76×35=2660", 2) | Out-Null
$d.Content.Find.Execute("84×12=1008", $true, $false, $false, $false, $false, $true, 1, $false, "91×66=6006", 2) | Out-Null
$d.Content.Find.Execute("56×61=3416", $true, $false, $false, $false, $false, $true, 1, $false, "66×57=3762", 2) | Out-Null
$d.Content.Find.Execute("17×78=1326", $true, $false, $false, $false, $false, $true, 1, $false, "63×75=4725", 2) | Out-Null
$d.Content.Find.Execute("60×59=3540", $true, $false, $false, $false, $false, $true, 1, $false, "33×95=3135", 2) | Out-Null
$d.Content.Find.Execute("25×60=1500", $true, $false, $false, $false, $false, $true, 1, $false, "72×32=2304", 2) | Out-Null
$d.Content.Find.Execute("57×56=3192", $true, $false, $false, $false, $false, $true, 1, $false, "14×16=224", 2) | Out-Null
$d.Content.Find.Execute("54×70=3780", $true, $false, $false, $false, $false, $true, 1, $false, "94×83=7802", 2) | Out-Null
$d.Content.Find.Execute("39×35=1365", $true, $false, $false, $false, $false, $true, 1, $false, "32×22=704", 2) | Out-Null
$d.Content.Find.Execute("82×41=3362", $true, $false, $false, $false, $false, $true, 1, $false, "78×41=3198", 2) | Out-Null
$d.Content.Find.Execute("69×64=4416", $true, $false, $false, $false, $false, $true, 1, $false, "24×49=1176", 2) | Out-Null
$d.Content.Find.Execute("41×20=820", $true, $false, $false, $false, $false, $true, 1, $false, "98×60=5880", 2) | Out-Null
$d.Content.Find.Execute("97×27=2619", $true, $false, $false, $false, $false, $true, 1, $false, "18×27=486", 2) | Out-Null
$d.Content.Find.Execute("46×24=1104", $true, $false, $false, $false, $false, $true, 1, $false, "48×50=2400", 2) | Out-Null
$d.Content.Find.Execute("65×58=3770", $true, $false, $false, $false, $false, $true, 1, $false, "73×21=1533", 2) | Out-Null
$d.Content.Find.Execute("34×94=3196", $true, $false, $false, $false, $false, $true, 1, $false, "73×95=6935", 2) | Out-Null
$d.Content.Find.Execute("81×83=6723", $true, $false, $false, $false, $false, $true, 1, $false, "98×20=1960", 2) | Out-Null
$d.Content.Find.Execute("59×49=2891", $true, $false, $false, $false, $false, $true, 1, $false, "43×74=3182", 2) | Out-Null
$d.Content.Find.Execute("11×88=968", $true, $false, $false, $false, $false, $true, 1, $false, "56×62=3472", 2) | Out-Null
$d.Content.Find.Execute("13×15=195", $true, $false, $false, $false, $false, $true, 1, $false, "97×25=2425", 2) | Out-Null
$d.Content.Find.Execute("42×81=3402", $true, $false, $false, $false, $false, $true, 1, $false, "60×50=3000", 2) | Out-Null
$d.Content.Find.Execute("40×26=1040", $true, $false, $false, $false, $false, $true, 1, $false, "53×16=848", 2) | Out-Null
$d.Content.Find.Execute("88×56=4928", $true, $false, $false, $false, $false, $true, 1, $false, "15×98=1470", 2) | Out-Null
